$d = $word.ActiveDocument

# The document originally contains two near-identical copies of the same
# sample content (a "Heading1" title followed by a body paragraph, a
# "Heading2" heading, another body paragraph, a page break, a second
# "Heading2" heading and a third body paragraph). The edit removes the
# first copy's paragraphs entirely (paragraphs 2-8), and replaces the
# text of the very first (Heading1) paragraph with the new title,
# leaving the second copy intact as the remaining body of the document.

$firstPara = $d.Paragraphs(1)
$eighthPara = $d.Paragraphs(8)

# Delete everything from just after the first paragraph's text through
# the end of the (former) eighth paragraph - i.e. paragraphs 2..8.
$deleteRange = $d.Range($firstPara.Range.End, $eighthPara.Range.End)
$deleteRange.Delete()

# Update the remaining first paragraph's text (keeps its Heading1 style).
$d.Paragraphs(1).Range.Text = "Zimbani Monthly sales report"
